# Apply cryptos list update (prices / 1h volume %) per commit diff.
# Values are written with a leading quote-prefix so Excel treats
# numeric-looking strings (e.g. "65.139.16", "596.85") as literal
# text instead of auto-converting them to numbers; ClearFormats()
# immediately after strips the quote-prefix/style side effect so the
# cell's style index is left untouched (matches original formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.139.16"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +3.39%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.630.50"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +2.17%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'596.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +1.76%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'155.45"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +5.27%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  -0.04%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  +1.19%  "
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'  +9.15%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  +5.76%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'5.78"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +1.37%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  +2.14%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'29.12"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +6.99%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.0000187"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +23.25%  "
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'3.102.77"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +2.25%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'64.994.06"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +3.35%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'2.634.07"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +2.54%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'12.53"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +3.58%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'4.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +3.90%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'352.05"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +2.72%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'7.35"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +8.84%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  +0.25%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'68.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +2.79%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'9.53"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +5.47%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'1.65"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -2.08%  "
$ws.Range("E25").ClearFormats()
$ws.Range("E26").Value = "'  +0.17%  "
$ws.Range("E26").ClearFormats()
$ws.Range("E27").Value = "'  +1.77%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'8.06"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +0.96%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'0.0₃0953"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +13.49%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  +0.01%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'528.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -4.22%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +4.99%  "
$ws.Range("E32").ClearFormats()
$ws.Range("E33").Value = "'  +2.92%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'5.55"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +8.46%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'6.32"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +7.14%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.426"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +4.14%  "
$ws.Range("E36").ClearFormats()
$ws.Range("B37").Value = "'EthereumClassic"
$ws.Range("B37").ClearFormats()
$ws.Range("C37").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C37").ClearFormats()
$ws.Range("D37").Value = "'20.31"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +5.35%  "
$ws.Range("E37").ClearFormats()
$ws.Range("B38").Value = "'Monero"
$ws.Range("B38").ClearFormats()
$ws.Range("C38").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C38").ClearFormats()
$ws.Range("D38").Value = "'163.85"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -0.83%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'1.99"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +6.39%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -0.02%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  -0.03%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'42.28"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +6.88%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'165.27"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.33%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'4.09"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  +4.45%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.0616"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +5.97%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'23.03"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +2.63%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'2.21"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +9.76%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.646"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +3.14%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.0255"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +3.81%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0980"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +2.35%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'19.41"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +3.08%  "
$ws.Range("E51").ClearFormats()
